$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E4").Value = 4
$ws.Range("C6").Select()
